# Add files via upload
# - Adds a new "Sheet2" after "Sheet1" with a 5-column header row
#   (IdOfOwner, LineNum, PackageType, StartDate, EndDate)
# - Makes Sheet2 the active/selected sheet (tabSelected moves off Sheet1)
# - Updates Sheet1's selection from A6 -> A2
# - Sets Sheet2's zoom to 160% and a couple of column widths
# - Leaves the selection on Sheet2 at D7

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New worksheet, placed right after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row
$ws2.Range("A1").Value = "IdOfOwner"
$ws2.Range("B1").Value = "LineNum"
$ws2.Range("C1").Value = "PackageType"
$ws2.Range("D1").Value = "StartDate"
$ws2.Range("E1").Value = "EndDate"

# Column widths (host quantizes ColumnWidth to 1/6-character steps, so these
# inputs are chosen to land as close as possible to the target 14.7109375 /
# 12.42578125 stored widths)
$ws2.Columns.Item(1).ColumnWidth = 13.833333333333334
$ws2.Columns.Item(3).ColumnWidth = 11.666666666666666

# Update the selection left behind on Sheet1
$ws1.Range("A2").Select()

# Activate Sheet2, zoom it, and leave its selection on D7
$ws2.Activate()
$ws2.Range("D7").Select()
$excel.ActiveWindow.Zoom = 160

Write-Host "Sheet2 added; workbook now has" $wb.Worksheets.Count "sheets"
